# Automatische test-sync: 2025-08-13 21:55:50
# Appends two new log rows (19 & 20) to the "Logs" sheet, extends the
# conditional-formatting ranges to cover them, and bumps the matching
# category counter on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# --- Row 19: fallback / manual-follow-up notification row -------------
$ws.Range("A19").Value = "[Fallback] Handmatige opvolging: Demo inplannen"
$ws.Range("B19").Value = "admin@testbedrijf123.nl"
$ws.Range("C19").Value = "Beste collega,`nOnderstaande e-mail kon niet automatisch worden beantwoord door het AI-systeem. Wil je deze even opvolgen?`n📩 Originele afzender: klantenservice@testbedrijf123.nl`n📝 Onderwerp: Demo inplannen`n🔎 Reden: Interne afzender`n━━━━━━━━━━━━━━━━━━━━━━━━━━━`n✉️ Bericht:`nKun je vrijdag om 11:00 een demo inplannen bij Van Dijk?`n━━━━━━━━━━━━━━━━━━━━━━━━━━━`nMet vriendelijke groet,`nMailMind Automatische Assistent`n—`n[Bedrijfsnaam]`nklantenservice@bedrijf.nl`nwww.bedrijf.nl"
$ws.Range("D19").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E19").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F19").Value = "2025-08-13 21:55:30"
$ws.Range("G19").Value = "Nee"
$ws.Range("H19").Value = "Ja"
$ws.Range("I19").Value = "Nee"
$ws.Range("J19").Value = "Nee"

# --- Row 20: original inbound request row ------------------------------
$ws.Range("A20").Value = "Demo inplannen"
$ws.Range("B20").Value = "klantenservice@testbedrijf123.nl"
$ws.Range("C20").Value = "Kun je vrijdag om 11:00 een demo inplannen bij Van Dijk?"
$ws.Range("D20").Value = "Intern verzoek / Actie voor medewerker"
$ws.Range("E20").Value = "Bedankt, we hebben dit doorgestuurd naar planning@testbedrijf123.nl."
$ws.Range("F20").Value = "2025-08-13 21:55:31"
$ws.Range("G20").Value = "Nee"
$ws.Range("H20").Value = "Ja"
$ws.Range("I20").Value = "Nee"
$ws.Range("J20").Value = "Nee"

# --- Extend conditional formatting ranges from row 18 to row 20 --------
$ws.Range("D2:D18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D20"))
$ws.Range("G2:G18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G20"))
$ws.Range("H2:H18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H20"))
$ws.Range("I2:I18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I20"))
$ws.Range("J2:J18").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J20"))

# --- Dashboard: bump the "Intern verzoek / Actie voor medewerker" count 
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 19

